$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.135.88"
$ws.Range("E2").Value = "  +0.69%  "
$ws.Range("D3").Value = "1.679.50"
$ws.Range("E3").Value = "  +0.35%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'215.09"
$ws.Range("E5").Value = "  +0.18%  "
$ws.Range("E6").Value = "  +0.20%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").Value = "'0.256"
$ws.Range("E8").Value = "  +2.07%  "
$ws.Range("D9").Value = "'21.48"
$ws.Range("E9").Value = "  +5.58%  "
$ws.Range("E10").Value = "  +0.73%  "
$ws.Range("D11").Value = "'0.0889"
$ws.Range("E11").Value = "  +0.23%  "
$ws.Range("D12").Value = "1.914.82"
$ws.Range("E12").Value = "  +0.27%  "
$ws.Range("D13").Value = "1.709.64"
$ws.Range("E13").Value = "  +1.99%  "
$ws.Range("E14").Value = "  +1.44%  "
$ws.Range("E15").Value = "  +1.91%  "
$ws.Range("D16").Value = "'66.29"
$ws.Range("E16").Value = "  +0.91%  "
$ws.Range("D17").Value = "27.117.63"
$ws.Range("E17").Value = "  +0.56%  "
$ws.Range("D18").Value = "'238.62"
$ws.Range("E18").Value = "  +0.86%  "
$ws.Range("E19").Value = "  +0.35%  "
$ws.Range("E20").Value = "  +1.35%  "
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("E22").Value = "  +2.24%  "
$ws.Range("E23").Value = "  +2.94%  "
$ws.Range("D24").Value = "'2.11"
$ws.Range("E24").Value = "  -3.44%  "
$ws.Range("D25").Value = "'147.83"
$ws.Range("E25").Value = "  +1.57%  "
$ws.Range("D26").Value = "'7.26"
$ws.Range("E26").Value = "  +0.33%  "
$ws.Range("D27").Value = "'16.30"
$ws.Range("E27").Value = "  +2.10%  "
$ws.Range("E28").Value = "  +0.41%  "
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("E30").Value = "  +0.16%  "
$ws.Range("E31").Value = "  +0.05%  "
$ws.Range("D32").Value = "1.568.16"
$ws.Range("E32").Value = "  +5.96%  "
$ws.Range("E33").Value = "  +1.60%  "
$ws.Range("E34").Value = "  +2.66%  "
$ws.Range("E35").Value = "  +0.51%  "
$ws.Range("E36").Value = "  +2.78%  "
$ws.Range("E37").Value = "  -1.17%  "
$ws.Range("D38").Value = "'0.934"
$ws.Range("E38").Value = "  +4.48%  "
$ws.Range("E39").Value = "  +1.04%  "
$ws.Range("E40").Value = "  +2.98%  "
$ws.Range("D41").Value = "'68.98"
$ws.Range("E41").Value = "  +3.06%  "
$ws.Range("E42").Value = "  -0.07%  "
$ws.Range("E43").Value = "  -5.00%  "
$ws.Range("D44").Value = "'2.25"
$ws.Range("E44").Value = "  -2.19%  "
$ws.Range("D45").Value = "1.823.75"
$ws.Range("E45").Value = "  +0.42%  "
$ws.Range("D46").Value = "'0.784"
$ws.Range("E46").Value = "  +1.18%  "
$ws.Range("D47").Value = "'90.71"
$ws.Range("E47").Value = "  +0.25%  "
$ws.Range("E48").Value = "  +3.00%  "
$ws.Range("E49").Value = "  +0.89%  "
$ws.Range("E50").Value = "  +5.87%  "
$ws.Range("E51").Value = "  +1.75%  "
